# Fruta / hortaliza, semanal
# Insert 4 new weekly rows of "Chirimoya" price data at the top of the
# data block (row 41), pushing the existing rows 41-72 down to 45-76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 41:72 down by 4 rows (-> 45:76), inserting 4 blank rows.
$ws.Range("A41:A44").EntireRow.Insert()

# Fill the 4 newly inserted rows with this week's data. All the
# "descriptive" columns (A,B,C,E,F,G,H,I,J,K,Q,R,T) are identical to every
# other row in this data block, only the date (D), quality (L), volume (M),
# prices (N,O,P) and price/kg (S) change per quality grade.

$newRows = @(
    @{ Row = 41; D = 44806; L = "Especial"; M = 150; N = 28000; S = 2800 },
    @{ Row = 42; D = 44806; L = "Primera";  M = 100; N = 25000; S = 2500 },
    @{ Row = 43; D = 44806; L = "Segunda";  M = 60;  N = 22000; S = 2200 },
    @{ Row = 44; D = 44806; L = "Tercera";  M = 30;  N = 18000; S = 1800 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.N
    $ws.Cells.Item($row, 16).Value = $r.N
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 10 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 10
}
